$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# ------------------------------------------------------------------
# 1) Apply formatting (borders / wrap-text) to the touched cells FIRST,
#    using PasteSpecial(formats) from existing template cells so the
#    styles reuse the workbook's existing border/alignment combinations.
#    This step never touches cell VALUES, so it has no effect on the
#    shared-string table ordering below.
# ------------------------------------------------------------------

# Row 17 & 18: B/C formatting normalizes to the plain bordered style
# (B: border only, C: border + wrap text), and a new blank E cell
# (border only) is added.
$ws.Range("A9").Copy() | Out-Null
$ws.Range("B17:B18").PasteSpecial(-4122) | Out-Null
$ws.Range("C4").Copy() | Out-Null
$ws.Range("C17:C18").PasteSpecial(-4122) | Out-Null
$ws.Range("A9").Copy() | Out-Null
$ws.Range("E17:E18").PasteSpecial(-4122) | Out-Null

# Rows 19-21: same column formatting pattern as the rest of the table
# (A: border only, B: border only, C: border + wrap, D: border only,
# E: border only).
$ws.Range("A9").Copy() | Out-Null
$ws.Range("A19:B21").PasteSpecial(-4122) | Out-Null
$ws.Range("C4").Copy() | Out-Null
$ws.Range("C19:C21").PasteSpecial(-4122) | Out-Null
$ws.Range("A9").Copy() | Out-Null
$ws.Range("D19:E21").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# ------------------------------------------------------------------
# 2) Populate cell VALUES. The order below matches the order the new
#    test-case data was authored in, which drives the order new
#    strings are appended to the shared-string table.
# ------------------------------------------------------------------

$ws.Range("C19").Value = "Verify that ""Add alternative name"" button should be disabled (Gryed`n * out) until a single letter is entered in last name field"
$ws.Range("A19").Value = "WAT18"
$ws.Range("B19").Value = "WAT-191"
$ws.Range("A20").Value = "WAT19"
$ws.Range("A21").Value = "WAT20"
$ws.Range("C20").Value = "Verify that ""Add alternative name"" button should be in disabled state even if First name field is having value."
$ws.Range("C21").Value = "Verify that ""Add alternative name"" button should be (Gryed out) when content of Last Name field is removed."
$ws.Range("B21").Value = "WAT-206"
$ws.Range("B20").Value = "WAT-180"

$ws.Range("D19").Value = "Y"
$ws.Range("D20").Value = "Y"
$ws.Range("D21").Value = "Y"

# Row 19's description contains a hard line-break, so it displays on two
# lines at the sheet's default row height.
$ws.Rows.Item(19).RowHeight = 30

# ------------------------------------------------------------------
# 3) Match the final selection left behind by the edit.
# ------------------------------------------------------------------
$ws.Range("A1:E21").Select()
